$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C (the "Förändrad" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# Increment every date in C2:C<lastRow> by one day (45179 -> 45180)
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value2
    if ($cur -ne $null) {
        $cell.Value2 = $cur + 1
    }
}
